$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Section 1: "Validation scores" -> "Validation accuracy scores" table
# ---------------------------------------------------------------------------

# Title
$ws.Range("A1").Value = "Validation accuracy scores"

# Table1 (Tableau3, A2:E10) header renames
$lo1 = $ws.ListObjects.Item(1)
$lo1.HeaderRowRange.Item(1,3).Value = "small (%)"
$lo1.HeaderRowRange.Item(1,4).Value = "medium (%)"
$lo1.HeaderRowRange.Item(1,5).Value = "large (%)"

# Fix model-variant label casing/text
$ws.Range("A3").Value = "LSTM Embed"
$ws.Range("A4").Value = "LST Vanilla"

# Clear the stray empty placeholder cells in column A for the sub-rows
$ws.Range("A6").Clear() | Out-Null
$ws.Range("A7").Clear() | Out-Null
$ws.Range("A9").Clear() | Out-Null
$ws.Range("A10").Clear() | Out-Null
$ws.Range("A16").Clear() | Out-Null
$ws.Range("A17").Clear() | Out-Null
$ws.Range("A19").Clear() | Out-Null
$ws.Range("A20").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Section 2: "Pretraining" table (Tableau6) - add medium/large columns + rows
# ---------------------------------------------------------------------------

$lo2 = $ws.ListObjects.Item(2)

# Resize table from A14:C20 to A14:E22 (adds 2 columns + 2 rows)
$lo2.Resize($ws.Range("A14:E22"))

# Header row
$lo2.HeaderRowRange.Item(1,3).Value = "small (%)"
$lo2.HeaderRowRange.Item(1,4).Value = "medium (%)"
$lo2.HeaderRowRange.Item(1,5).Value = "large (%)"

# Row 15 (LSTM Embed) - extend body highlight style to the new D/E cells
$ws.Range("C15").Copy($ws.Range("D15:E15")) | Out-Null

# Row 16 (overlap window slide) - new pretraining score formula + N/A cells
$ws.Range("C16").Formula = "=(0.0332 + 0.0421 + 0.0391) / 3 * 100"
$ws.Range("C16").NumberFormat = "#,##0.00_);(#,##0.00)"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "N/A"

# Row 17 (full window slide) - new pretraining score formula + N/A cells
$ws.Range("C17").Formula = "=(0.0307 + 0.019 + 0.0146) / 3 * 100"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "N/A"

# Row 18 (LSTM Vanilla) - extend body highlight style to the new D/E cells
$ws.Range("C18").Copy($ws.Range("D18:E18")) | Out-Null

# Row 19/20 - N/A cells
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = "N/A"
$ws.Range("D20").Value = "N/A"
$ws.Range("E20").Value = "N/A"

# New rows 21/22 - Embedding / Embedding augmented
# Copy row 3's style pattern (A=body-highlight, B=body-highlight+border, C/D/E=body-highlight)
$ws.Range("A3:E3").Copy($ws.Range("A21:E21")) | Out-Null
$ws.Range("A21").Value = "Embedding"
$ws.Range("B21").ClearContents() | Out-Null
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "N/A"
$ws.Range("E21").ClearContents() | Out-Null

$ws.Range("A3:E3").Copy($ws.Range("A22:E22")) | Out-Null
$ws.Range("A22").Value = "Embedding augmented"
$ws.Range("B22").ClearContents() | Out-Null
$ws.Range("C22").Value = "N/A"
$ws.Range("D22").Value = "N/A"
$ws.Range("E22").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Selection marker to mirror final author state
# ---------------------------------------------------------------------------
$ws.Range("M18").Select() | Out-Null
